$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 157, shifting existing rows 157..177 down to 158..178
$ws.Rows.Item(157).Insert()

# Populate the new row 157 with the weekly Acelga record
$ws.Cells.Item(157, 1).Value = 4
$ws.Cells.Item(157, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(157, 3).Value = "Los Lagos"
$ws.Cells.Item(157, 4).Value = 44711
$ws.Cells.Item(157, 5).Value = 10
$ws.Cells.Item(157, 6).Value = 100112009
$ws.Cells.Item(157, 7).Value = "Acelga"
$ws.Cells.Item(157, 8).Value = "Sin especificar"
$ws.Cells.Item(157, 9).Value = "Primera"
$ws.Cells.Item(157, 10).Value = 40
$ws.Cells.Item(157, 11).Value = 12000
$ws.Cells.Item(157, 12).Value = 12000
$ws.Cells.Item(157, 13).Value = 12000
$ws.Cells.Item(157, 14).Value = '$/docena de atados (12 kilos)'
$ws.Cells.Item(157, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(157, 16).Value = 1000
$ws.Cells.Item(157, 17).Value = 12
$ws.Cells.Item(157, 18).Value = "Hortaliza"
